# Apply the edits described in the commit:
# "moved subsets "HeatSlowRamper", "HeatQuickRamper", "Hydro", and "Biomass" to file.
#  deleted subsets "Offshore", "Onshore", and "Oil", since not used"
#
# This replaces the placeholder "Lignite" values that filled column B on sheet
# "Par_TagTechnologyToSubsets" for rows 469-866 with the real subset names, and
# populates column A (Technology) for the first 37 of those rows with the
# corresponding technology names (taken from quoted GAMS set text, mirroring the
# existing rows above them).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Par_TagTechnologyToSubsets")
$ws2 = $wb.Worksheets.Item("Par_TagFuelToSubsets")

# ---------------------------------------------------------------------------
# Technology names (column A) for rows 469-505, grouped by destination subset
# ---------------------------------------------------------------------------
$heatSlowRamperTechs = @(
    "HLR_Oil_Boiler",
    "HLI_Oil_Boiler",
    "HHI_BF_BOF",
    "HHI_DRI_EAF",
    "HHI_Scrap_EAF",
    "HHI_H2DRI_EAF",
    "HHI_Molten_Electrolysis",
    "HHI_Bio_BF_BOF",
    "HHI_BF_BOF_CCS",
    "HHI_DRI_EAF_CCS"
)

$heatQuickRamperTechs = @(
    "HLR_Hardcoal",
    "HLR_Lignite",
    "HLR_Biomass",
    "HLR_Gas_Boiler",
    "HLR_Direct_Electric",
    "HLR_H2_Boiler",
    "HLI_Hardcoal",
    "HLI_Lignite",
    "HLI_Biomass",
    "HLI_Gas_Boiler",
    "HLI_Direct_Electric",
    "HLI_H2_Boiler",
    "HMI_Gas",
    "HMI_Steam_Electric",
    "HMI_Gas_CCS",
    "HMI_Biomass",
    "HMI_HardCoal",
    "HMI_Oil",
    "HMI_HardCoal_CCS"
)

$hydroTechs = @(
    "Res_Hydro_large",
    "Res_Hydro_small"
)

$biomassTechsWithName = @(
    "RES_Grass",
    "RES_Wood",
    "RES_Residues",
    "RES_Paper_Cardboard",
    "RES_Roundwood",
    "RES_Biogas"
)

# Row where each block starts, and the subset name (column B) for the block
$startRow = 469
$blocks = @(
    @{ Techs = $heatSlowRamperTechs;  Subset = "HeatSlowRamper" },
    @{ Techs = $heatQuickRamperTechs; Subset = "HeatQuickRamper" },
    @{ Techs = $hydroTechs;           Subset = "Hydro" },
    @{ Techs = $biomassTechsWithName; Subset = "Biomass" }
)

# Style-source cell that already carries the Consolas / vertical-center style
# used for every populated "Technology" cell in column A.
$styleSource = $ws1.Range("A454")

$row = $startRow
foreach ($block in $blocks) {
    foreach ($tech in $block.Techs) {
        $styleSource.Copy() | Out-Null
        $ws1.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null
        $ws1.Cells.Item($row, 1).Value = $tech
        $ws1.Cells.Item($row, 2).Value = $block.Subset
        $row++
    }
}

# ---------------------------------------------------------------------------
# Remaining rows (506-866) only get their Subset (column B) updated to
# "Biomass" - column A stays empty, exactly as it was before the edit.
# ---------------------------------------------------------------------------
for ($r = 506; $r -le 866; $r++) {
    $ws1.Cells.Item($r, 2).Value = "Biomass"
}

# ---------------------------------------------------------------------------
# Sheet view / selection bookkeeping to mirror the authored file
# ---------------------------------------------------------------------------
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 477
$ws1.Range("H497").Select() | Out-Null

$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 3

$ws1.Activate()
